$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to hold literal text (no numeric/date auto-conversion),
# then restore the "Normal" style so no stray number-format sticks to it.
function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "66.666.36"
$ws.Range("E2").Value = "  -4.33%  "
Set-TextValue "D3" "3.338.80"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "573.73"
$ws.Range("E5").Value = "  -3.38%  "
Set-TextValue "D6" "180.70"
$ws.Range("E6").Value = "  -5.54%  "
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("E10").Value = "  -1.60%  "
Set-TextValue "D11" "0.405"
$ws.Range("E11").Value = "  -3.56%  "
Set-TextValue "D12" "3.916.56"
$ws.Range("E12").Value = "  -1.36%  "
Set-TextValue "D13" "0.135"
$ws.Range("E13").Value = "  -0.63%  "
Set-TextValue "D14" "27.04"
$ws.Range("E14").Value = "  -5.53%  "
Set-TextValue "D15" "66.762.41"
$ws.Range("E15").Value = "  -4.21%  "
$ws.Range("E16").Value = "  -2.82%  "
Set-TextValue "D17" "3.337.71"
$ws.Range("E17").Value = "  -1.68%  "
Set-TextValue "D18" "437.46"
$ws.Range("E18").Value = "  -3.59%  "
$ws.Range("E19").Value = "  -2.26%  "
Set-TextValue "D20" "13.55"
$ws.Range("E20").Value = "  -2.01%  "
Set-TextValue "D21" "7.59"
$ws.Range("E21").Value = "  -2.56%  "
Set-TextValue "D22" "73.53"
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("E23").Value = "  +0.12%  "
Set-TextValue "D24" "0.517"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").Value = "  -4.14%  "
$ws.Range("E26").Value = "  -0.18%  "
Set-TextValue "D27" "9.03"
$ws.Range("E27").Value = "  -5.08%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("E31").Value = "  +0.05%  "
Set-TextValue "D32" "5.27"
$ws.Range("E32").Value = "  -5.96%  "
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("E34").Value = "  -4.68%  "
Set-TextValue "D35" "163.83"
$ws.Range("E35").Value = "  -0.42%  "
Set-TextValue "D36" "1.48"
$ws.Range("E36").Value = "  -5.95%  "
Set-TextValue "D37" "27.37"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("E38").Value = "  -8.44%  "
Set-TextValue "D39" "2.834.15"
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("E40").Value = "  -2.11%  "
Set-TextValue "D41" "4.42"
$ws.Range("E41").Value = "  -3.87%  "
Set-TextValue "D42" "6.22"
$ws.Range("E42").Value = "  -5.86%  "
Set-TextValue "D43" "40.21"
$ws.Range("E43").Value = "  -2.19%  "
Set-TextValue "D44" "0.0667"
$ws.Range("E44").Value = "  -3.03%  "
Set-TextValue "D45" "24.41"
$ws.Range("E45").Value = "  -4.26%  "
$ws.Range("E46").Value = "  -7.03%  "
Set-TextValue "D47" "321.94"
$ws.Range("E47").Value = "  -5.11%  "
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D50" "0.976"
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D51" "6.17"
$ws.Range("E51").Value = "  -2.64%  "
